# SCANClim_Configuration.xlsx update
# - Pest_list: swap example pest from "Fusarium odoratissimum" to "Agrilus planipennis"
# - Other settings: "Remove climates not in EU" flipped from "yes" to "no"
# - Climates_to_be_removed: drop the extra climate rows (Ocean, Dsb, Dsc), keeping just the header
# - misc selection / active-sheet bookkeeping to match the saved workbook view state

$wb = $excel.ActiveWorkbook

# --- Pest_list: update the example pest name ---
$pestList = $wb.Worksheets.Item("Pest_list")
$pestList.Range("A2").Value = "Agrilus planipennis"

# --- Other settings: "Remove climates not in EU" -> "no" ---
$otherSettings = $wb.Worksheets.Item("Other settings")
$otherSettings.Range("B2").Value = "no"

# --- Climates_to_be_removed: remove the Ocean/Dsb/Dsc rows, leaving only the "Climates" header ---
$climates = $wb.Worksheets.Item("Climates_to_be_removed")
$climates.Rows("2:4").Delete()

# --- Restore per-sheet selections to match the latest saved view state ---
$authors = $wb.Worksheets.Item("Authors")
$authors.Range("A2:A3").Select()

$pestStatus = $wb.Worksheets.Item("Pest_status_to_be_included")
$pestStatus.Range("G34").Select()

$otherSettings.Range("B13").Select()

$pestList.Range("A2").Select()

# Climates_to_be_removed becomes the active/selected sheet last, matching tabSelected + activeTab
$climates.Range("A2:A4").Select()
